$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.713252999999999
$ws.Range("H2").Value = 20.139759
$ws.Range("I2").Value = 0.3101840064655811
$ws.Range("J2").Value = 0.3231642354899327
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07487166666666667
$ws.Range("N2").Value = 0.224615
$ws.Range("O2").Value = 0.01287435003490057
$ws.Range("P2").Value = 0.01655871537719798
$ws.Range("Q2").Value = 0.5026324408649999
$ws.Range("R2").Value = 4.523691967785
$ws.Range("S2").Value = 0.003993417474465753
$ws.Range("T2").Value = 0.005351184595567576
# Row 3
$ws.Range("G3").Value = 6.713252999999999
$ws.Range("H3").Value = 20.139759
$ws.Range("I3").Value = 0.3101840064655811
$ws.Range("J3").Value = 0.3231642354899327
$ws.Range("O3").Value = 0.02600892111095355
$ws.Range("P3").Value = 0.03345212152666174
$ws.Range("Q3").Value = 1.01542427127
$ws.Range("R3").Value = 9.138818441429999
$ws.Range("S3").Value = 0.008067551354042806
$ws.Range("T3").Value = 0.01081052927867996
# Row 4
$ws.Range("G4").Value = 6.713252999999999
$ws.Range("H4").Value = 20.139759
$ws.Range("I4").Value = 0.3101840064655811
$ws.Range("J4").Value = 0.3231642354899327
$ws.Range("M4").Value = 0.8000470000000001
$ws.Range("N4").Value = 2.400141
$ws.Range("O4").Value = 0.1375698656239178
$ws.Range("P4").Value = 0.1769394371887155
$ws.Range("Q4").Value = 5.370917922891
$ws.Range("R4").Value = 48.338261306019
$ws.Range("S4").Value = 0.04267197208815843
$ws.Range("T4").Value = 0.05718049794711019
# Row 5
$ws.Range("G5").Value = 6.713252999999999
$ws.Range("H5").Value = 20.139759
$ws.Range("I5").Value = 0.3101840064655811
$ws.Range("J5").Value = 0.3231642354899327
$ws.Range("M5").Value = 3.881946
$ws.Range("N5").Value = 7.763892
$ws.Range("O5").Value = 0.6675092708044715
$ws.Range("P5").Value = 0.5723574910282232
$ws.Range("Q5").Value = 26.060485630338
$ws.Range("R5").Value = 156.362913782028
$ws.Range("S5").Value = 0.2070506999710496
$ws.Range("T5").Value = 0.1849654710150717
# Row 6
$ws.Range("G6").Value = 6.713252999999999
$ws.Range("H6").Value = 20.139759
$ws.Range("I6").Value = 0.3101840064655811
$ws.Range("J6").Value = 0.3231642354899327
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9074473333333333
$ws.Range("N6").Value = 2.722342
$ws.Range("O6").Value = 0.1560375924257564
$ws.Range("P6").Value = 0.2006922348792017
$ws.Range("Q6").Value = 6.091923532841999
$ws.Range("R6").Value = 54.82731179557799
$ws.Range("S6").Value = 0.04840036557786455
$ws.Range("T6").Value = 0.06485655265350321
# Row 7
$ws.Range("I7").Value = 0.0154484264788496
$ws.Range("J7").Value = 0.01609489473505086
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07487166666666667
$ws.Range("N7").Value = 0.224615
$ws.Range("O7").Value = 0.01287435003490057
$ws.Range("P7").Value = 0.01655871537719798
$ws.Range("Q7").Value = 0.02503314209222222
$ws.Range("R7").Value = 0.22529827883
$ws.Range("S7").Value = 0.0001988884499771362
$ws.Range("T7").Value = 0.0002665107809436693
# Row 8
$ws.Range("I8").Value = 0.0154484264788496
$ws.Range("J8").Value = 0.01609489473505086
$ws.Range("O8").Value = 0.02600892111095355
$ws.Range("P8").Value = 0.03345212152666174
$ws.Range("Q8").Value = 0.05057226314888889
$ws.Range("R8").Value = 0.45515036834
$ws.Range("S8").Value = 0.0004017969055767651
$ws.Range("T8").Value = 0.0005384083746357494
# Row 9
$ws.Range("I9").Value = 0.0154484264788496
$ws.Range("J9").Value = 0.01609489473505086
$ws.Range("M9").Value = 0.8000470000000001
$ws.Range("N9").Value = 2.400141
$ws.Range("O9").Value = 0.1375698656239178
$ws.Range("P9").Value = 0.1769394371887155
$ws.Range("Q9").Value = 0.2674935809913334
$ws.Range("R9").Value = 2.407442228922
$ws.Range("S9").Value = 0.002125237954796312
$ws.Range("T9").Value = 0.002847821616031518
# Row 10
$ws.Range("I10").Value = 0.0154484264788496
$ws.Range("J10").Value = 0.01609489473505086
$ws.Range("M10").Value = 3.881946
$ws.Range("N10").Value = 7.763892
$ws.Range("O10").Value = 0.6675092708044715
$ws.Range("P10").Value = 0.5723574910282232
$ws.Range("Q10").Value = 1.297918293244
$ws.Range("R10").Value = 7.787509759464
$ws.Range("S10").Value = 0.01031196789397338
$ws.Range("T10").Value = 0.009212033568917068
# Row 11
$ws.Range("I11").Value = 0.0154484264788496
$ws.Range("J11").Value = 0.01609489473505086
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.9074473333333333
$ws.Range("N11").Value = 2.722342
$ws.Range("O11").Value = 0.1560375924257564
$ws.Range("P11").Value = 0.2006922348792017
$ws.Range("Q11").Value = 0.3034025960404444
$ws.Range("R11").Value = 2.730623364364
$ws.Range("S11").Value = 0.002410535274525997
$ws.Range("T11").Value = 0.003230120394522853
# Row 12
$ws.Range("G12").Value = 6.661784666666667
$ws.Range("H12").Value = 19.985354
$ws.Range("I12").Value = 0.3078059262949933
$ws.Range("J12").Value = 0.3206866401135023
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.07487166666666667
$ws.Range("N12").Value = 0.224615
$ws.Range("O12").Value = 0.01287435003490057
$ws.Range("P12").Value = 0.01655871537719798
$ws.Range("Q12").Value = 0.4987789209677778
$ws.Range("R12").Value = 4.48901028871
$ws.Range("S12").Value = 0.003962801237938549
$ws.Range("T12").Value = 0.005310158798909404
# Row 13
$ws.Range("G13").Value = 6.661784666666667
$ws.Range("H13").Value = 19.985354
$ws.Range("I13").Value = 0.3078059262949933
$ws.Range("J13").Value = 0.3206866401135023
$ws.Range("O13").Value = 0.02600892111095355
$ws.Range("P13").Value = 0.03345212152666174
$ws.Range("Q13").Value = 1.007639342731111
$ws.Range("R13").Value = 9.06875408458
$ws.Range("S13").Value = 0.008005700054490463
$ws.Range("T13").Value = 0.01072764845705372
# Row 14
$ws.Range("G14").Value = 6.661784666666667
$ws.Range("H14").Value = 19.985354
$ws.Range("I14").Value = 0.3078059262949933
$ws.Range("J14").Value = 0.3206866401135023
$ws.Range("M14").Value = 0.8000470000000001
$ws.Range("N14").Value = 2.400141
$ws.Range("O14").Value = 0.1375698656239178
$ws.Range("P14").Value = 0.1769394371887155
$ws.Range("Q14").Value = 5.329740837212667
$ws.Range("R14").Value = 47.96766753491401
$ws.Range("S14").Value = 0.04234481991864777
$ws.Range("T14").Value = 0.05674211361562324
# Row 15
$ws.Range("G15").Value = 6.661784666666667
$ws.Range("H15").Value = 19.985354
$ws.Range("I15").Value = 0.3078059262949933
$ws.Range("J15").Value = 0.3206866401135023
$ws.Range("M15").Value = 3.881946
$ws.Range("N15").Value = 7.763892
$ws.Range("O15").Value = 0.6675092708044715
$ws.Range("P15").Value = 0.5723574910282232
$ws.Range("Q15").Value = 25.860688339628
$ws.Range("R15").Value = 155.164130037768
$ws.Range("S15").Value = 0.2054633094104659
$ws.Range("T15").Value = 0.1835474007416349
# Row 16
$ws.Range("G16").Value = 6.661784666666667
$ws.Range("H16").Value = 19.985354
$ws.Range("I16").Value = 0.3078059262949933
$ws.Range("J16").Value = 0.3206866401135023
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.9074473333333333
$ws.Range("N16").Value = 2.722342
$ws.Range("O16").Value = 0.1560375924257564
$ws.Range("P16").Value = 0.2006922348792017
$ws.Range("Q16").Value = 6.045218731007555
$ws.Range("R16").Value = 54.406968579068
$ws.Range("S16").Value = 0.04802929567345059
$ws.Range("T16").Value = 0.06435931850028102
# Row 17
$ws.Range("G17").Value = 2.607918
$ws.Range("H17").Value = 5.215835999999999
$ws.Range("I17").Value = 0.1204981331366039
$ws.Range("J17").Value = 0.08369373503331734
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.07487166666666667
$ws.Range("N17").Value = 0.224615
$ws.Range("O17").Value = 0.01287435003490057
$ws.Range("P17").Value = 0.01655871537719798
$ws.Range("Q17").Value = 0.19525916719
$ws.Range("R17").Value = 1.17155500314
$ws.Range("S17").Value = 0.00155133514455269
$ws.Range("T17").Value = 0.001385860737271325
# Row 18
$ws.Range("G18").Value = 2.607918
$ws.Range("H18").Value = 5.215835999999999
$ws.Range("I18").Value = 0.1204981331366039
$ws.Range("J18").Value = 0.08369373503331734
$ws.Range("O18").Value = 0.02600892111095355
$ws.Range("P18").Value = 0.03345212152666174
$ws.Range("Q18").Value = 0.39446498362
$ws.Range("R18").Value = 2.36678990172
$ws.Range("S18").Value = 0.003134026438767108
$ws.Range("T18").Value = 0.002799732995354758
# Row 19
$ws.Range("G19").Value = 2.607918
$ws.Range("H19").Value = 5.215835999999999
$ws.Range("I19").Value = 0.1204981331366039
$ws.Range("J19").Value = 0.08369373503331734
$ws.Range("M19").Value = 0.8000470000000001
$ws.Range("N19").Value = 2.400141
$ws.Range("O19").Value = 0.1375698656239178
$ws.Range("P19").Value = 0.1769394371887155
$ws.Range("Q19").Value = 2.086456972146
$ws.Range("R19").Value = 12.518741832876
$ws.Range("S19").Value = 0.01657691198353555
$ws.Range("T19").Value = 0.01480872237301665
# Row 20
$ws.Range("G20").Value = 2.607918
$ws.Range("H20").Value = 5.215835999999999
$ws.Range("I20").Value = 0.1204981331366039
$ws.Range("J20").Value = 0.08369373503331734
$ws.Range("M20").Value = 3.881946
$ws.Range("N20").Value = 7.763892
$ws.Range("O20").Value = 0.6675092708044715
$ws.Range("P20").Value = 0.5723574910282232
$ws.Range("Q20").Value = 10.123796848428
$ws.Range("R20").Value = 40.49518739371199
$ws.Range("S20").Value = 0.08043362098331459
$ws.Range("T20").Value = 0.04790273619845042
# Row 21
$ws.Range("G21").Value = 2.607918
$ws.Range("H21").Value = 5.215835999999999
$ws.Range("I21").Value = 0.1204981331366039
$ws.Range("J21").Value = 0.08369373503331734
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.9074473333333333
$ws.Range("N21").Value = 2.722342
$ws.Range("O21").Value = 0.1560375924257564
$ws.Range("P21").Value = 0.2006922348792017
$ws.Range("Q21").Value = 2.366548234651999
$ws.Range("R21").Value = 14.199289407912
$ws.Range("S21").Value = 0.01880223858643393
$ws.Range("T21").Value = 0.0167966827292242
# Row 22
$ws.Range("G22").Value = 5.325505333333333
$ws.Range("H22").Value = 15.976516
$ws.Range("I22").Value = 0.2460635076239721
$ws.Range("J22").Value = 0.2563604946281968
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.07487166666666667
$ws.Range("N22").Value = 0.224615
$ws.Range("O22").Value = 0.01287435003490057
$ws.Range("P22").Value = 0.01655871537719798
$ws.Range("Q22").Value = 0.3987294601488889
$ws.Range("R22").Value = 3.58856514134
$ws.Range("S22").Value = 0.003167907727966442
$ws.Range("T22").Value = 0.004245000464506001
# Row 23
$ws.Range("G23").Value = 5.325505333333333
$ws.Range("H23").Value = 15.976516
$ws.Range("I23").Value = 0.2460635076239721
$ws.Range("J23").Value = 0.2563604946281968
$ws.Range("O23").Value = 0.02600892111095355
$ws.Range("P23").Value = 0.03345212152666174
$ws.Range("Q23").Value = 0.8055181850355556
$ws.Range("R23").Value = 7.24966366532
$ws.Range("S23").Value = 0.006399846358076408
$ws.Range("T23").Value = 0.008575802420937551
# Row 24
$ws.Range("G24").Value = 5.325505333333333
$ws.Range("H24").Value = 15.976516
$ws.Range("I24").Value = 0.2460635076239721
$ws.Range("J24").Value = 0.2563604946281968
$ws.Range("M24").Value = 0.8000470000000001
$ws.Range("N24").Value = 2.400141
$ws.Range("O24").Value = 0.1375698656239178
$ws.Range("P24").Value = 0.1769394371887155
$ws.Range("Q24").Value = 4.260654565417333
$ws.Range("R24").Value = 38.34589108875601
$ws.Range("S24").Value = 0.03385092367877971
$ws.Range("T24").Value = 0.04536028163693385
# Row 25
$ws.Range("G25").Value = 5.325505333333333
$ws.Range("H25").Value = 15.976516
$ws.Range("I25").Value = 0.2460635076239721
$ws.Range("J25").Value = 0.2563604946281968
$ws.Range("M25").Value = 3.881946
$ws.Range("N25").Value = 7.763892
$ws.Range("O25").Value = 0.6675092708044715
$ws.Range("P25").Value = 0.5723574910282232
$ws.Range("Q25").Value = 20.673324126712
$ws.Range("R25").Value = 124.039944760272
$ws.Range("S25").Value = 0.1642496725456681
$ws.Range("T25").Value = 0.146729849504149
# Row 26
$ws.Range("G26").Value = 5.325505333333333
$ws.Range("H26").Value = 15.976516
$ws.Range("I26").Value = 0.2460635076239721
$ws.Range("J26").Value = 0.2563604946281968
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.9074473333333333
$ws.Range("N26").Value = 2.722342
$ws.Range("O26").Value = 0.1560375924257564
$ws.Range("P26").Value = 0.2006922348792017
$ws.Range("Q26").Value = 4.832615613385777
$ws.Range("R26").Value = 43.493540520472
$ws.Range("S26").Value = 0.03839515731348137
$ws.Range("T26").Value = 0.06485655265350321

Write-Host "Applied all changes"